$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9347826086956522
$ws.Range("C3").Value = 43
$ws.Range("D3").Value = 43
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 3
$ws.Range("J3").Value = "wonderful"
$ws.Range("K3").Value = 0.8392857142857143
$ws.Range("L3").Value = 47
$ws.Range("M3").Value = 47
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 9
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.8636363636363636
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 38
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 6
$ws.Range("J4").Value = "awesome"
$ws.Range("K4").Value = 0.8153846153846154
$ws.Range("L4").Value = 53
$ws.Range("M4").Value = 53
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 12
$ws.Range("A5").Value = "however"
$ws.Range("B5").Value = 0.75
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 48
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 16
$ws.Range("J5").Value = "favorite"
$ws.Range("K5").Value = 0.6989247311827957
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 28
$ws.Range("A6").Value = "disappointed"
$ws.Range("B6").Value = 0.7311827956989247
$ws.Range("C6").Value = 136
$ws.Range("D6").Value = 136
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = "classic"
$ws.Range("K6").Value = 0.6226415094339622
$ws.Range("L6").Value = 33
$ws.Range("M6").Value = 33
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 20
$ws.Range("A7").Value = "broke"
$ws.Range("B7").Value = 0.7281553398058253
$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 150
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 56
$ws.Range("J7").Value = "excellent"
$ws.Range("K7").Value = 0.53125
$ws.Range("L7").Value = 34
$ws.Range("M7").Value = 34
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 30
$ws.Range("A8").Value = "poor"
$ws.Range("B8").Value = 0.6901408450704225
$ws.Range("C8").Value = 49
$ws.Range("D8").Value = 49
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 22
$ws.Range("J8").Value = "great"
$ws.Range("K8").Value = 0.340983606557377
$ws.Range("L8").Value = 416
$ws.Range("M8").Value = 416
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 804
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.6216216216216216
$ws.Range("C9").Value = 92
$ws.Range("D9").Value = 92
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 56
$ws.Range("J9").Value = "love"
$ws.Range("K9").Value = 0.2801724137931034
$ws.Range("L9").Value = 195
$ws.Range("M9").Value = 196
$ws.Range("N9").Value = 0.99
$ws.Range("O9").Value = 0.01000000000000001
$ws.Range("P9").Value = $true
$ws.Range("Q9").Value = 501
$ws.Range("A10").Value = "smaller"
$ws.Range("B10").Value = 0.5546218487394958
$ws.Range("C10").Value = 66
$ws.Range("D10").Value = 66
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 53
$ws.Range("J10").Value = "loves"
$ws.Range("K10").Value = 0.2572614107883817
$ws.Range("L10").Value = 124
$ws.Range("M10").Value = 124
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 358
$ws.Range("A11").Value = "junk"
$ws.Range("B11").Value = 0.5272727272727272
$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 29
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 26
$ws.Range("J11").Value = "loved"
$ws.Range("K11").Value = 0.1896024464831804
$ws.Range("L11").Value = 62
$ws.Range("M11").Value = 62
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 265
$ws.Range("A12").Value = "broken"
$ws.Range("B12").Value = 0.4939759036144578
$ws.Range("C12").Value = 41
$ws.Range("D12").Value = 41
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 42
$ws.Range("J12").Value = "perfect"
$ws.Range("K12").Value = 0.1807228915662651
$ws.Range("L12").Value = 30
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 136
$ws.Range("A13").Value = "small"
$ws.Range("B13").Value = 0.472463768115942
$ws.Range("C13").Value = 163
$ws.Range("D13").Value = 163
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 182
$ws.Range("J13").Value = "friends"
$ws.Range("K13").Value = 0.164021164021164
$ws.Range("L13").Value = 31
$ws.Range("M13").Value = 31
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 158
$ws.Range("A14").Value = "apart"
$ws.Range("B14").Value = 0.4210526315789473
$ws.Range("C14").Value = 40
$ws.Range("D14").Value = 40
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 55
$ws.Range("J14").Value = "fun"
$ws.Range("K14").Value = 0.08764241893076249
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 1041
$ws.Range("A15").Value = "plastic"
$ws.Range("B15").Value = 0.4094488188976378
$ws.Range("C15").Value = 52
$ws.Range("D15").Value = 52
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 75
$ws.Range("J15").Value = "game"
$ws.Range("K15").Value = 0.03634003893575601
$ws.Range("L15").Value = 56
$ws.Range("M15").Value = 56
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 1485
$ws.Range("A16").Value = "thought"
$ws.Range("B16").Value = 0.301980198019802
$ws.Range("C16").Value = 61
$ws.Range("D16").Value = 61
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 141
$ws.Range("A17").Value = "ok"
$ws.Range("B17").Value = 0.28125
$ws.Range("C17").Value = 36
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 92
$ws.Range("A18").Value = "cheap"
$ws.Range("B18").Value = 0.2796208530805687
$ws.Range("C18").Value = 59
$ws.Range("D18").Value = 59
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 152
$ws.Range("A19").Value = "size"
$ws.Range("B19").Value = 0.2371134020618557
$ws.Range("C19").Value = 46
$ws.Range("D19").Value = 46
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 148
$ws.Range("A20").Value = "item"
$ws.Range("B20").Value = 0.2028985507246377
$ws.Range("C20").Value = 56
$ws.Range("D20").Value = 56
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 220
$ws.Range("A21").Value = "hard"
$ws.Range("B21").Value = 0.195
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 161
$ws.Range("A22").Value = "money"
$ws.Range("B22").Value = 0.1930379746835443
$ws.Range("C22").Value = 61
$ws.Range("D22").Value = 61
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 255
$ws.Range("A23").Value = "would"
$ws.Range("B23").Value = 0.1824925816023739
$ws.Range("C23").Value = 123
$ws.Range("D23").Value = 123
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 551
$ws.Range("A24").Value = "used"
$ws.Range("B24").Value = 0.1714285714285714
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 145
$ws.Range("A25").Value = "work"
$ws.Range("B25").Value = 0.1582278481012658
$ws.Range("C25").Value = 50
$ws.Range("D25").Value = 50
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 266
$ws.Range("A26").Value = "product"
$ws.Range("B26").Value = 0.1387665198237885
$ws.Range("C26").Value = 63
$ws.Range("D26").Value = 63
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 391
$ws.Range("A27").Value = "price"
$ws.Range("B27").Value = 0.117816091954023
$ws.Range("C27").Value = 41
$ws.Range("D27").Value = 41
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 307
$ws.Range("A28").Value = "2"
$ws.Range("B28").Value = 0.1161048689138577
$ws.Range("C28").Value = 31
$ws.Range("D28").Value = 31
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 236
$ws.Range("A29").Value = "like"
$ws.Range("B29").Value = 0.05766062602965404
$ws.Range("C29").Value = 35
$ws.Range("D29").Value = 36
$ws.Range("E29").Value = 0.03
$ws.Range("F29").Value = 0.97
$ws.Range("G29").Value = $true
$ws.Range("H29").Value = 572


# Row 16 no longer has the second block of data (J:Q) - clear those cells
$ws.Range("J16:Q16").ClearContents()

# Apply the header-row style (bold, centered, bordered) to the new A29 label cell
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
